# Added Simple Icons icon
# Insert a new row 37 holding the "Simple Icons" entry, which pushes the
# previously-existing rows 37-42 down to 38-43 (keeping their data intact).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 37.
$ws.Rows("37:37").Insert()

# Populate the new row with the Simple Icons entry.
$ws.Range("A37").Value = "Simple Icons"
$ws.Range("B37").Value = "BF1813"
$ws.Range("C37").Value = 2
$ws.Range("D37").Formula = "=MOD((C37+100),360)"
$ws.Range("E37").Value = 90
$ws.Range("F37").Value = 75

# Move the selection the same way the source workbook shows after the edit.
$ws.Range("A44").Select() | Out-Null
